$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42 : new journal entry (06.06.2018) ---
$ws.Range("A42").Value = 43256
$ws.Range("B42").Value = 6
$ws.Range("C42").Value = "Modification de l'interface"
$ws.Range("D42").Value = "L'interface donne maintenant (en temps réel) la largeur et la hauteur de la selectionbox, et les labels deviennent verts quand les proportions de l'échantillon sélectionné permettent d'enregistrer celui-ci sans (trop) le déformer. J'ai aussi commencé à développer un système permettant à l'utilisateur de lire les caractères, et enregistrer les valeurs lues dans un fichier csv, permettant ainsi de faire des paires échantillon/label, qui serviront ensuite à l'entraînement de réseau."
$ws.Range("E42").Value = 6

# match the formatting used on the row above (same visual "block" style)
$ws.Range("A37").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("B37").Copy()
$ws.Range("B42").PasteSpecial(-4122)
$ws.Range("C37").Copy()
$ws.Range("C42").PasteSpecial(-4122)
$ws.Range("D37").Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E42").PasteSpecial(-4122)

$ws.Rows.Item(42).RowHeight = 90

# --- Row 43 : new journal entry (07.06.2018) ---
$ws.Range("A43").Value = 43257
$ws.Range("B43").Value = 2
$ws.Range("C43").Value = "Développement d'un système d'importation de datasets"
$ws.Range("D43").Value = "L'application est désormais capable d'importer un dataset composé d'un dossier rempli d'images, et d'un fichier csv"

$ws.Rows.Item(43).RowHeight = 30

$excel.CutCopyMode = 0

# --- sheet view : move to the newly-edited area ---
$ws.Application.ActiveWindow.ScrollRow = 34
[void]$ws.Range("D44").Select()

Write-Output "done"
